$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.19057559967041
$ws.Range("B1").Value = 1.895533084869385
$ws.Range("C1").Value = 6.401522636413574
$ws.Range("D1").Value = 2.284320592880249
$ws.Range("E1").Value = 1.193129420280457
